$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.662.13'
$ws.Range("E2").Value = '  -1.72%  '

$ws.Range("D3").Value = '3.329.24'
$ws.Range("E3").Value = '  -3.84%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '617.75'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.06%  '

$ws.Range("E7").Value = '  -4.89%  '

$ws.Range("E8").Value = '  -3.48%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.927'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.34%  '

$ws.Range("D11").Value = '3.328.65'
$ws.Range("E11").Value = '  -3.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.14'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '

$ws.Range("E13").Value = '  -2.43%  '

$ws.Range("D14").Value = '93.497.10'
$ws.Range("E14").Value = '  -1.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.95'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.67%  '

$ws.Range("D16").Value = '3.946.53'
$ws.Range("E16").Value = '  -3.99%  '

$ws.Range("E17").Value = '  -4.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.47%  '

$ws.Range("D19").Value = '3.320.64'
$ws.Range("E19").Value = '  -3.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +8.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '494.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.78%  '

$ws.Range("E24").Value = '  -12.99%  '

$ws.Range("E25").Value = '  -5.33%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.40'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.25%  '

$ws.Range("E28").Value = '  -2.82%  '

$ws.Range("D29").Value = '3.505.06'
$ws.Range("E29").Value = '  -3.71%  '

$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  +1.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.82%  '

$ws.Range("E33").Value = '  -4.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.174'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -8.97%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.529'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '533.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.60%  '

$ws.Range("E39").Value = '  +0.09%  '

$ws.Range("E40").Value = '  -4.95%  '

$ws.Range("E41").Value = '  -1.89%  '

$ws.Range("E42").Value = '  -5.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.862'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.34%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.74%  '

$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.29%  '

$ws.Range("E49").Value = '  -2.25%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.05'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.95'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.02%  '
